$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '249.05'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '22.72'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.275'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05704'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.406'
$ws.Range("B7").Value = 'KuCoinToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '6.335'
$ws.Range("E7").Value = '6KuCoinTokenKCS'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8052'
$ws.Range("E8").Value = '7MXTokenMX'
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9008'
$ws.Range("E9").Value = '8FTXTokenFTT'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1425'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07457'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03092'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03001'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09388'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.860'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001575'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04802'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("B18").Value = 'UpBots'
$ws.Range("C18").Value = 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.01827'
$ws.Range("E18").Value = '17UpBotsUBXTBestin24h'
$ws.Range("B19").Value = 'One'
$ws.Range("C19").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0005810'
$ws.Range("E19").Value = '18OneONE'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.006431'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0009976'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.0001499'
$ws.Range("B25").Value = 'BTSEToken'
$ws.Range("C25").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.201'
$ws.Range("E25").Value = '24BTSETokenBTSE'
$ws.Range("B26").Value = 'BitpandaEcosystemToken'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.3300'
$ws.Range("E26").Value = '25BitpandaEcosystemTokenBEST'
$ws.Range("B27").Value = 'ProBitToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1292'
$ws.Range("E27").Value = '26ProBitTokenPROB'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03982'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006827'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1070'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002729'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.007699'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005565'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.2069'
